# Add a new worksheet "L6" (summary of last-6-games Form / Goals scored /
# Goals conceded / Total Goals per team) after the last existing sheet.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "L6"

# Header row
$ws.Range("B1").Value = "Form"
$ws.Range("C1").Value = "Goals scored"
$ws.Range("D1").Value = "Goals conceded"
$ws.Range("E1").Value = "Total Goals"

$data = @(
    @("1","Belenenses,L D W D W L","Belenenses,0 1 2 2 2 0","Belenenses,2 1 0 2 1 1","Belenenses,2 2 2 4 3 1"),
    @("2","Benfica,W W L W W W","Benfica,1 5 1 5 2 2","Benfica,0 0 2 1 1 0","Benfica,1 5 3 6 3 2"),
    @("3","Boavista,W D W L L D","Boavista,2 3 2 1 0 3","Boavista,0 3 0 2 1 3","Boavista,2 6 2 3 1 6"),
    @("4","Famalicao,W D L W D L","Famalicao,2 1 0 3 2 2","Famalicao,0 1 1 0 2 3","Famalicao,2 2 1 3 4 5"),
    @("5","Farense,L L L W D D","Farense,1 0 0 2 1 0","Farense,2 1 1 0 1 0","Farense,3 1 1 2 2 0"),
    @("6","Gil Vicente,W L W L L D","Gil Vicente,2 1 2 0 1 0","Gil Vicente,0 2 1 3 2 0","Gil Vicente,2 3 3 3 3 0"),
    @("7","Guimaraes,L L W L L W","Guimaraes,1 0 1 0 0 2","Guimaraes,2 3 0 1 1 0","Guimaraes,3 3 1 1 1 2"),
    @("8","Maritimo,L W L W W W","Maritimo,0 1 0 1 1 1","Maritimo,1 0 2 0 0 0","Maritimo,1 1 2 1 1 1"),
    @("9","Moreirense,D W L D D L","Moreirense,1 2 2 0 1 0","Moreirense,1 1 3 0 1 2","Moreirense,2 3 5 0 2 2"),
    @("10","Nacional,L L L L W L","Nacional,1 1 0 1 1 0","Nacional,5 5 1 2 0 2","Nacional,6 6 1 3 1 2"),
    @("11","Pacos Ferreira,L L L L D W","Pacos Ferreira,0 0 0 0 1 1","Pacos Ferreira,2 5 2 2 1 0","Pacos Ferreira,2 5 2 2 2 1"),
    @("12","Portimonense,W W W L D D","Portimonense,5 3 1 1 1 0","Portimonense,1 0 0 5 1 0","Portimonense,6 3 1 6 2 0"),
    @("13","Porto,W W W W D W","Porto,2 2 1 1 1 3","Porto,1 0 0 0 1 2","Porto,3 2 1 1 2 5"),
    @("14","Rio Ave,L D D L D D","Rio Ave,0 3 0 0 1 0","Rio Ave,2 3 0 1 1 0","Rio Ave,2 6 0 1 2 0"),
    @("15","Santa Clara,L W L D L D","Santa Clara,1 5 0 0 1 3","Santa Clara,2 1 1 0 2 3","Santa Clara,3 6 1 0 3 6"),
    @("16","Sp Braga,W D D W L L","Sp Braga,2 1 0 2 0 0","Sp Braga,1 1 0 1 1 1","Sp Braga,3 2 0 3 1 1"),
    @("17","Sp Lisbon,D D W D W W","Sp Lisbon,1 1 1 2 1 2","Sp Lisbon,1 1 0 2 0 0","Sp Lisbon,2 2 1 4 1 2"),
    @("18","Tondela,W L W W D L","Tondela,2 0 3 2 2 0","Tondela,1 2 2 1 2 2","Tondela,3 2 5 3 4 2")
)

$row = 2
foreach ($r in $data) {
    # Column A holds the rank number, stored as text (matches the rest of
    # the workbook, e.g. the "Table" sheet's A column).
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $r[0]

    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]

    $row = $row + 1
}
